# Change the letters to lowercase
# The "blok_"/"piek_" name labels in column B of Sheet1 (rows 2-19) contain
# segments "GGG", "GLG" and "GHG" that should be lowercased to "ggg", "glg"
# and "ghg" respectively, while leaving the rest of the text (including the
# "T10"/"T100"/"T1000" suffix) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $value = $cell.Value2
    if ($null -ne $value) {
        $newValue = $value -creplace "GGG", "ggg" -creplace "GLG", "glg" -creplace "GHG", "ghg"
        if (-not $newValue.Equals($value)) {
            $cell.Value2 = $newValue
        }
    }
}
